$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13: "Longest subarray with given sum K(positives)" -> Website=GFG, Level=Medium
$ws.Range("C13").Value = "GFG"
$ws.Range("D13").Value = "Medium"

# Row 14: "Longest subarray with sum K (Positives + Negatives)" -> Website=GFG, Level=Medium
$ws.Range("C14").Value = "GFG"
$ws.Range("D14").Value = "Medium"

# Row 15: "Longest Consecutive Sequence in an Array" -> Website=LC, Level=Medium
$ws.Range("C15").Value = "LC"
$ws.Range("D15").Value = "Medium"

# Update the active selection to D15 to match the saved cursor position
$ws.Range("D15").Select()
